# test data for one sample ID is substring of another sample ID
#
# Renames the sample pair that previously was:
#   test_sample_2_T_IGO / test_investigator_sample_2_T / test_sample_2_T
#   test_sample_1_N_IGO / test_investigator_sample_1_N / test_sample_1_N
# to a pair where one id is a substring of the other:
#   test_sample_1a_IGO / test_investigator_sample_1a / test_sample_1a
#   test_sample_1_IGO  / test_investigator_sample_1  / test_sample_1

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # SampleInfo
$ws2 = $wb.Worksheets.Item(2)   # SampleRenames

# --- SampleInfo sheet: update the CMO sample id / investigator sample id cells
# (write investigator-sample columns first, then CMO-sample-id columns, to
# reproduce the shared-string insertion order of the original edit)
$ws1.Range("C2").Value = "test_investigator_sample_1a"
$ws1.Range("C3").Value = "test_investigator_sample_1"
$ws1.Range("A2").Value = "test_sample_1a_IGO"
$ws1.Range("A3").Value = "test_sample_1_IGO"

# --- SampleRenames sheet: update the OldName / NewName cells
$ws2.Range("A2").Value = "test_sample_1a_IGO"
$ws2.Range("A3").Value = "test_sample_1_IGO"
$ws2.Range("B2").Value = "test_sample_1a"
$ws2.Range("B3").Value = "test_sample_1"

# --- View state: the active tab moves from SampleInfo to SampleRenames,
# and the selections/scroll position on each sheet change.
$ws1.Range("A3").Select()
$ws2.Activate()
$ws2.Range("B13").Select()
